$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.654
$ws.Range("C5").Value = -12.156
$ws.Range("C9").Value = -11.671
$ws.Range("C11").Value = -12.704
$ws.Range("A21").Value = -21.785
$ws.Range("C21").Value = -13.336
$ws.Range("A23").Value = -21.654
$ws.Range("A25").Value = -21.937
